# Add a new fiscal-year (FY2018, period ending 2018-12-31) column to the
# CALX yearly financials sheet. This is implemented the same way a user
# would do it in Excel: insert a new column D (pushing the existing D:K
# data one column to the right, into E:L) and then fill the new column D
# with the FY2018 figures for Income Statement, Balance Sheet and Cash
# Flow Statement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at D; everything from D:K shifts to E:L.
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D cells default to the generic style.
#    Copy the number/date formatting from column E (the column that used
#    to be D before the insert) onto column D for every data row, so the
#    new column keeps the same look (date format row 7/38/80, #,##0
#    format everywhere else) instead of falling back to "General".
$ws.Range("E7:E35").Copy($ws.Range("D7:D35"))
$ws.Range("E38:E77").Copy($ws.Range("D38:D77"))
$ws.Range("E80:E102").Copy($ws.Range("D80:D102"))

# 3) Populate the new column D with the FY2018 values.

# --- Income Statement ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 441300
$ws.Range("D9").Value = 243900
$ws.Range("D10").Value = 197400
$ws.Range("D12").Value = 90000
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 5700
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 459800
$ws.Range("D18").Value = -18500
$ws.Range("D20").Value = -300
$ws.Range("D21").Value = -9600
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = -18800
$ws.Range("D24").Value = 500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -19300
$ws.Range("D27").Value = -19300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 300
$ws.Range("D33").Value = -19300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -19300

# --- Balance Sheet ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 45800
$ws.Range("D42").Value = 3800
$ws.Range("D43").Value = 67000
$ws.Range("D44").Value = 50200
$ws.Range("D45").Value = 7900
$ws.Range("D46").Value = 174800
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 24900
$ws.Range("D49").Value = 116200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 317100
$ws.Range("D57").Value = 40200
$ws.Range("D58").Value = 30000
$ws.Range("D59").Value = 73500
$ws.Range("D60").Value = 143700
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 21500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D66").Value = 165100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -684900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 151900
$ws.Range("D77").Value = 0

# --- Cash Flow Statement ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -19300
$ws.Range("D83").Value = 9200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 3600
$ws.Range("D91").Value = -10400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -100
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 7500
$ws.Range("D101").Value = -500
$ws.Range("D102").Value = 10500
